# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 15, pushing the existing
# rows 15-46 down to 16-47 (Excel's normal Insert-row shift behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..46 down to 16..47, leaving a blank row 15 in place.
$ws.Rows("15:15").Insert()

# Populate the new row 15 with the new weekly record (same market /
# product / unit metadata as its neighbours; new date + price figures).
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 45260
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112039
$ws.Range("G15").Value = "Ciboulette"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2500
$ws.Range("N15").Value = "$/docena de atados"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 833
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = "Hortaliza"

Write-Output "Inserted row 15; sheet now spans to row 47."
